$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 10, shifting existing rows 10+ down by one.
$ws.Rows.Item(10).Insert() | Out-Null

# Fill the new row 10 with the "In-fill gain lifetime" data.
$ws.Range("A10").Value = "In-fill gain lifetime"
$ws.Range("B10").Value = 5
$ws.Range("C10").Value = 11.2
$ws.Range("D10").Value = 11.6
$ws.Range("E10").Value = 16.5

$ws.Range("G20").Select() | Out-Null
